$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Section heading: "Visual Check:" -> "Performance Verification:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Visual Check:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Performance Verification:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Re-word the existing "Visual Check" table rows (Table 3)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("No damage to cuff or tubing", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Visual Inspection", 2) | Out-Null
$d.Content.Find.Execute("Velcro integrity", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Functional Check", 2) | Out-Null
$d.Content.Find.Execute("Check overpressure cut-off", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Quantity Checked", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "Sphygmo Test" table (Table 4): drop the first row (0mmhg check) and
#    reword the remaining two rows.
# ---------------------------------------------------------------------------
$tSphygmo = $d.Tables.Item(4)
$tSphygmo.Rows.Item(1).Delete()

$d.Content.Find.Execute("Sphygmo with zero press. Reads 100mmhg", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Ball Moving Freely, Not Sticking", 2) | Out-Null
$d.Content.Find.Execute("Sphygmo with zero press. Reads 200mmhg", $true, $false, $false, $false, $false, `
    $true, 1, $false, "No Leaks, Occlude O2 Nipple, Flow Drops to Zero", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "Leak Test" table (Table 5): reword its single row.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Leak Test no Significant Drop in Press @200mmhg", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Accuracy ± 5%", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Renumber the result placeholders that used to belong to the two tables
#    being folded into the first one (one result field was removed along
#    with the deleted 0mmhg row, so everything after it shifts down by one).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("<result5>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<result4>", 2) | Out-Null
$d.Content.Find.Execute("<result6>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<result5>", 2) | Out-Null
$d.Content.Find.Execute("<result7>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<result6>", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Copy the (now re-worded) rows of tables 4 and 5 onto the end of table 3,
#    then get rid of the now-redundant tables and their heading paragraphs so
#    everything lives in a single "Performance Verification" table.
# ---------------------------------------------------------------------------
$tMain = $d.Tables.Item(3)

$rowsToCopy = @(
    @("Ball Moving Freely, Not Sticking", "<result4>"),
    @("No Leaks, Occlude O2 Nipple, Flow Drops to Zero", "<result5>"),
    @("Accuracy ± 5%", "<result6>"),
    @("O2 Nipple intact, No Damage", "<result7>"),
    @("Replaced Outlet O Ring", "<result8>"),
    @("Replaced Seal", "<result9>")
)

foreach ($pair in $rowsToCopy) {
    $newRow = $tMain.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $pair[0]
    $newRow.Cells.Item(2).Range.Text = $pair[1]
}

# Remove the now-empty "Sphygmo Test" and "Leak Test" tables outright.
$d.Tables.Item(5).Delete()
$d.Tables.Item(4).Delete()

# Remove the two section headings that used to introduce those tables
# (full paragraph, including its end-of-paragraph mark).
$paras = @()
foreach ($p in $d.Paragraphs) { $paras += $p }
for ($i = $paras.Count - 1; $i -ge 0; $i--) {
    $t = $paras[$i].Range.Text
    if (($t -like "*Sphygmo*Test:*") -or ($t -like "*Leak Test:*")) {
        $startPos = $paras[$i].Range.Start
        $endPos = $paras[$i].Range.End + 1
        $d.Range($startPos, $endPos).Delete()
    }
}

Write-Output "done"
